$d = $word.ActiveDocument

# --- Edit 1: add a new bullet point right under "REPLICATION OF THESIS GRAPHS:" ---
# Locate the heading paragraph, then take the paragraph that immediately
# follows it (the first existing bullet in that list) and insert a brand-new
# paragraph *before* it. Word copies that paragraph's mark formatting
# (pStyle=ListParagraph, numPr numId=2, bold paragraph-mark rPr) onto the
# freshly inserted paragraph, matching the other bullet items already in the
# list, exactly like the new bullet should look.
$headingIdx = -1
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "REPLICATION OF THESIS GRAPHS:") {
        $headingIdx = $i
        break
    }
}

$firstBulletIdx = $headingIdx + 1
$firstBullet = $d.Paragraphs.Item($firstBulletIdx)
$firstBullet.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($firstBulletIdx)
$newPara.Range.InsertAfter("Replaced original boxplots with normal error bars to reduce visual complexity and increase readability.")

# --- Edit 2: insert " journals" after "scientific" in the contrast-sensitivity sentence ---
$findRange = $d.Content
$found = $findRange.Find.Execute("scientific", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRange.Collapse(0)
    $findRange.InsertAfter(" journals")
}
